# Update view-count figures (column F) across sheets, per upstream data refresh.
$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F11").Value = 4797
$wsExpo.Range("F12").Value = 1362
$wsExpo.Range("F17").Value = 1190
$wsExpo.Range("F18").Value = 3946
$wsExpo.Range("F19").Value = 911
$wsExpo.Range("F35").Value = 130
$wsExpo.Range("F37").Value = 2059
$wsExpo.Range("F38").Value = 981
$wsExpo.Range("F41").Value = 550
$wsExpo.Range("F42").Value = 175

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F12").Value = 136

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 4797
$wsAll.Range("F17").Value = 1190
$wsAll.Range("F18").Value = 3946
$wsAll.Range("F19").Value = 911
$wsAll.Range("F36").Value = 2059
$wsAll.Range("F38").Value = 981
$wsAll.Range("F44").Value = 550
